# Auto-generated edit script: updates crypto price/volume table cells
# to match the target revision described in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a
# number (single-dot decimal-looking strings like "1.001" or "6.410",
# where e.g. a trailing zero would be silently dropped). For these we
# temporarily force a text number-format, assign the value, then put
# the cell style back to "Normal" so no visible formatting change sticks.
function Set-TextValue($rangeAddr, $value) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "22.531.50"
$ws.Range("E2").Value = "  +0.23%  "

# Row 3
$ws.Range("D3").Value = "1.574.78"
$ws.Range("E3").Value = "  +0.18%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
Set-TextValue "D5" "1.001"
$ws.Range("E5").Value = "  -0.01%  "

# Row 6
Set-TextValue "D6" "287.07"
$ws.Range("E6").Value = "  -1.61%  "

# Row 7
Set-TextValue "D7" "0.3665"
$ws.Range("E7").Value = "  -1.40%  "

# Row 8
Set-TextValue "D8" "48.46"
$ws.Range("E8").Value = "  -3.06%  "

# Row 9
Set-TextValue "D9" "0.3347"
$ws.Range("E9").Value = "  -1.43%  "

# Row 10
Set-TextValue "D10" "1.135"
$ws.Range("E10").Value = "  -0.76%  "

# Row 11
Set-TextValue "D11" "0.07456"
$ws.Range("E11").Value = "  -1.24%  "

# Row 12
Set-TextValue "D12" "1.001"
$ws.Range("E12").Value = "  -0.04%  "

# Row 13
Set-TextValue "D13" "20.91"
$ws.Range("E13").Value = "  -1.83%  "

# Row 14
Set-TextValue "D14" "6.006"
$ws.Range("E14").Value = "  -0.62%  "

# Row 15
Set-TextValue "D15" "6.938"
$ws.Range("E15").Value = "  -0.34%  "

# Row 16
$ws.Range("D16").Value = "1.578.14"
$ws.Range("E16").Value = "  +0.43%  "

# Row 17
Set-TextValue "D17" "0.00001114"
$ws.Range("E17").Value = "  -0.77%  "

# Row 18
Set-TextValue "D18" "88.46"
$ws.Range("E18").Value = "  -2.45%  "

# Row 19
Set-TextValue "D19" "0.06755"
$ws.Range("E19").Value = "  -0.10%  "

# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D20" "6.410"
$ws.Range("E20").Value = "  +1.81%  "

# Row 21
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D21" "1.002"
$ws.Range("E21").Value = "  +0.05%  "

# Row 22
Set-TextValue "D22" "16.49"
$ws.Range("E22").Value = "  +0.63%  "

# Row 23
$ws.Range("E23").Value = "  -0.22%  "

# Row 24
$ws.Range("D24").Value = "22.526.35"
$ws.Range("E24").Value = "  +0.22%  "

# Row 25
Set-TextValue "D25" "2.388"
$ws.Range("E25").Value = "  +1.13%  "

# Row 26
Set-TextValue "D26" "2.628"
$ws.Range("E26").Value = "  +0.09%  "

# Row 27
Set-TextValue "D27" "152.38"
$ws.Range("E27").Value = "  +1.91%  "

# Row 28
Set-TextValue "D28" "19.65"
$ws.Range("E28").Value = "  -2.02%  "

# Row 29
Set-TextValue "D29" "5.039"
$ws.Range("E29").Value = "  -0.36%  "

# Row 30
Set-TextValue "D30" "124.19"
$ws.Range("E30").Value = "  -0.82%  "

# Row 31
$ws.Range("D31").Value = "1.756.59"
$ws.Range("E31").Value = "  +0.52%  "

# Row 32
$ws.Range("E32").Value = "  -1.85%  "

# Row 33
Set-TextValue "D33" "6.204"
$ws.Range("E33").Value = "  -0.70%  "

# Row 34
Set-TextValue "D34" "2.002"
$ws.Range("E34").Value = "  -0.44%  "

# Row 35
Set-TextValue "D35" "9.828"
$ws.Range("E35").Value = "  +0.57%  "

# Row 36
Set-TextValue "D36" "0.08299"
$ws.Range("E36").Value = "  -0.77%  "

# Row 37
Set-TextValue "D37" "0.02459"
$ws.Range("E37").Value = "  -0.92%  "

# Row 38
Set-TextValue "D38" "0.2274"
$ws.Range("E38").Value = "  -1.28%  "

# Row 39
Set-TextValue "D39" "0.06462"
$ws.Range("E39").Value = "  -1.03%  "

# Row 40
Set-TextValue "D40" "5.466"
$ws.Range("E40").Value = "  +0.09%  "

# Row 41
Set-TextValue "D41" "1.301"
$ws.Range("E41").Value = "  -2.77%  "

# Row 42
Set-TextValue "D42" "0.6366"
$ws.Range("E42").Value = "  +2.08%  "

# Row 43
Set-TextValue "D43" "11.39"
$ws.Range("E43").Value = "  +0.39%  "

# Row 44
Set-TextValue "D44" "13.92"
$ws.Range("E44").Value = "  -0.59%  "

# Row 45
Set-TextValue "D45" "0.6228"
$ws.Range("E45").Value = "  +6.09%  "

# Row 46
Set-TextValue "D46" "3.763"
$ws.Range("E46").Value = "  -1.36%  "

# Row 47
Set-TextValue "D47" "2.062"
$ws.Range("E47").Value = "  -0.56%  "

# Row 48
Set-TextValue "D48" "125.35"
$ws.Range("E48").Value = "  -3.71%  "

# Row 49
Set-TextValue "D49" "1.222"
$ws.Range("E49").Value = "  +0.34%  "

# Row 50
Set-TextValue "D50" "0.07260"
$ws.Range("E50").Value = "  -0.96%  "

# Row 51
Set-TextValue "D51" "76.97"
$ws.Range("E51").Value = "  +0.39%  "

Write-Output "Applied crypto table update"
